# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff
# (Sun Jun  2 15:45:50 UTC 2024 GitHub Actions run).
#
# Every touched cell is forced to Text number-format before the
# assignment so values like "601.49" / "68.131.08" / "  +0.50%  "
# stay literal strings instead of being auto-coerced to numbers by
# Excel's usual text-looks-like-a-number heuristic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.131.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.795.77'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.49'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.14'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.84'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.431.63'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.780.03'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '68.114.78'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.37'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.45%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '461.47'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.703'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.04'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.02'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.944.42'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.36%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.24'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.35'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.35'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.04'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0997'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.33'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.85'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.987'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '47.62'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.22'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '152.52'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.36'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '392.90'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.61'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.06%  '
